$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Update the Version and Date property values (rows 3 and 8) ---
$ws1.Cells.Item(3, 2).Value = "0.2.0"
$ws1.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"

# --- Insert a new "Jurisdiction" row right after the "Contact" row (row 10), ---
# --- pushing "Description" and everything below it down by one row.         ---
$ws1.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (border/alignment) by
# copying the format from the row just below (which still carries the
# original style after the shift).
$ws1.Range("A12:B12").Copy() | Out-Null
$ws1.Range("A11:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = "iso:code:3166:FR"
